$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New "Date Watched" column (H) - header + per-row mixed-type values
#    (mirrors the "mixed data types" theme already used in columns A/B)
# ---------------------------------------------------------------------

$ws.Range("H1").Value = "Date Watched"

# Row 2: plain date, formatted as a date (mm-dd-yy -> numFmtId 14)
$ws.Range("H2").Value = 36194
$ws.Range("H2").NumberFormat = "mm-dd-yy"

# Row 3: numeric serial stored under a text format (stays numeric, not a string)
$ws.Range("H3").Value = 39160
$ws.Range("H3").NumberFormat = "@"

# Row 4: date
$ws.Range("H4").Value = 44330
$ws.Range("H4").NumberFormat = "mm-dd-yy"

# Row 5: literal text date (kept as text, not re-parsed into a serial number)
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "2013-01-25"

# Row 6: date
$ws.Range("H6").Value = 36892
$ws.Range("H6").NumberFormat = "mm-dd-yy"

# Row 7: numeric serial under a 2-decimal number format
$ws.Range("H7").Value = 40768
$ws.Range("H7").NumberFormat = "0.00"

# Row 8: date
$ws.Range("H8").Value = 39132
$ws.Range("H8").NumberFormat = "mm-dd-yy"

# Row 9: numeric serial under the General format
$ws.Range("H9").Value = 43071
$ws.Range("H9").Style = "Normal"

# Row 10: date
$ws.Range("H10").Value = 35749
$ws.Range("H10").NumberFormat = "mm-dd-yy"

# Column H width (auto-fit-ish)
$ws.Range("H1").EntireColumn.ColumnWidth = 12.71

# ---------------------------------------------------------------------
# 2. Column B ("Year Watched") - vary the cell formats across rows to
#    extend the mixed-data-types coverage (same underlying values).
# ---------------------------------------------------------------------

$ws.Range("B2").Value = 1996
$ws.Range("B2").NumberFormat = "0.00"
$ws.Range("B3").Value = 1975
$ws.Range("B3").NumberFormat = "0.00"
$ws.Range("B4").Value = 1956
$ws.Range("B4").NumberFormat = "0.00"

$ws.Range("B5").Value = 2007
$ws.Range("B5").Style = "Normal"
$ws.Range("B6").Value = 1968
$ws.Range("B6").Style = "Normal"
$ws.Range("B7").Value = 2009
$ws.Range("B7").Style = "Normal"

$ws.Range("B8").Value = 2006
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B9").Value = 2007
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B10").Value = 1937
$ws.Range("B10").NumberFormat = "@"

# ---------------------------------------------------------------------
# 3. Selection cosmetics
# ---------------------------------------------------------------------
$ws.Range("F16").Select()

Write-Output "done"
